$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.522.59'
$ws.Range('E2').Value = '  +3.71%  '
$ws.Range('D3').Value = '1.602.00'
$ws.Range('E3').Value = '  +3.00%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.15'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('E6').Value = '  +7.25%  '
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.83'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +10.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.52'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('E10').Value = '  +2.48%  '
$ws.Range('E11').Value = '  +2.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0912'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.23%  '
$ws.Range('D13').Value = '1.830.37'
$ws.Range('E13').Value = '  +2.94%  '
$ws.Range('D14').Value = '1.605.26'
$ws.Range('E14').Value = '  +3.19%  '
$ws.Range('D15').Value = '29.512.50'
$ws.Range('E15').Value = '  +3.66%  '
$ws.Range('E16').Value = '  +4.71%  '
$ws.Range('E17').Value = '  +3.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.54'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.32'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +4.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.60'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +3.15%  '
$ws.Range('E21').Value = '  +3.13%  '
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('E23').Value = '  +3.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.19'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.98%  '
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.43'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.30'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.52%  '
$ws.Range('E28').Value = '  +4.89%  '
$ws.Range('E29').Value = '  +2.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.998'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0472'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.43%  '
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('E33').Value = '  +2.51%  '
$ws.Range('D34').Value = '1.428.35'
$ws.Range('E34').Value = '  +2.44%  '
$ws.Range('E35').Value = '  +3.23%  '
$ws.Range('E36').Value = '  -1.15%  '
$ws.Range('E37').Value = '  +1.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.80'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +5.39%  '
$ws.Range('E40').Value = '  +1.99%  '
$ws.Range('E41').Value = '  +3.21%  '
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '53.91'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +23.64%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.795'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.34%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.997'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0472'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.79%  '
$ws.Range('E47').Value = '  +1.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.28'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.93%  '
$ws.Range('D49').Value = '1.740.89'
$ws.Range('E49').Value = '  +2.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '86.50'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.34%  '
$ws.Range('E51').Value = '  -3.64%  '
